$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where the existing "Java File Name" value in column D simply needs to
# swap places with the existing "Cpp File Name" value already present in
# column E (header row + the 5 questions that already had cpp code).
$swapRows = 1,3,4,5,6,7
foreach ($r in $swapRows) {
    $dVal = $ws.Cells.Item($r, 4).Value2
    $eVal = $ws.Cells.Item($r, 5).Value2
    $ws.Cells.Item($r, 4).Value = $eVal
    $ws.Cells.Item($r, 5).Value = $dVal
}

# Rows 8-14: new cpp files were added for questions 6-12. The java file name
# that used to live in column D moves over to column E, and the new cpp file
# name takes its place in column D.
$newCpp = @{
    8  = "UnionAndIntersection.cpp"
    9  = "RotateBy1.cpp"
    10 = "LargestSumContiguousSubarray.cpp"
    11 = "MinimizeTheHeights.cpp"
    12 = "MinNumOfJumpsToReachEnd.cpp"
    13 = "FindDuplicate.cpp"
    14 = "MergeSortedWithoutExtraSpace.cpp"
}
foreach ($r in 8..14) {
    $dVal = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($r, 5).Value = $dVal
    $ws.Cells.Item($r, 4).Value = $newCpp[$r]
}

# Row 15 duplicates row 10's entry (question 13 reuses the same cpp file
# reference as the "Largest sum contiguous subarray" question above it).
$row15Old = $ws.Cells.Item(15, 4).Value2
$ws.Cells.Item(15, 5).Value = $row15Old
$ws.Cells.Item(15, 4).Value = "LargestSumContiguousSubarray.cpp"

# Rows 16-18: no cpp file exists yet for these questions, so the java file
# name just moves from column D to column E, leaving D blank.
foreach ($r in 16..18) {
    $dVal = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($r, 5).Value = $dVal
    $ws.Cells.Item($r, 4).Clear()
}

# Column D/E widths swap along with the content swap.
$dWidth = $ws.Columns.Item(4).ColumnWidth
$eWidth = $ws.Columns.Item(5).ColumnWidth
$ws.Columns.Item(4).ColumnWidth = $eWidth
$ws.Columns.Item(5).ColumnWidth = $dWidth

# Update the current selection / view to match the edited area.
[void]$ws.Range("D15").Select()
